$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2099.6
$ws.Range("I43").Value = 1999.5
$ws.Range("J43").Value = 2166.3333
$ws.Range("K43").Value = 1999.5
$ws.Range("L43").Value = 2166.3333
$ws.Range("M43").Value = -1930.5
$ws.Range("N43").Value = -2304.3333
$ws.Range("H99").Value = 293.4762
$ws.Range("I99").Value = 258.15
$ws.Range("K99").Value = 774.4499999999999
$ws.Range("M99").Value = 723.5500000000001
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 251.16667
$ws.Range("J5").Value = 89
$ws.Range("L5").Value = 89
$ws.Range("N5").Value = -313
$ws.Range("H32").Value = 3400.923
$ws.Range("I32").Value = 2727.2083
$ws.Range("K32").Value = 2727.2083
$ws.Range("M32").Value = -2440.2083
$ws.Range("H40").Value = 25009.334
$ws.Range("I40").Value = 22514
$ws.Range("J40").Value = 30000
$ws.Range("K40").Value = 22514
$ws.Range("L40").Value = 30000
$ws.Range("M40").Value = -22338
$ws.Range("N40").Value = -30352
$ws.Range("H61").Value = 2692.3684
$ws.Range("I61").Value = 1929.5834
$ws.Range("K61").Value = 1929.5834
$ws.Range("M61").Value = -1717.5834
$ws.Range("H110").Value = 3745.7222
$ws.Range("J110").Value = 2222
$ws.Range("L110").Value = 2222
$ws.Range("N110").Value = -6312
$ws.Range("H132").Value = 2146.0789
$ws.Range("I132").Value = 1698.2307
$ws.Range("K132").Value = 5094.6921
$ws.Range("M132").Value = -2564.6921
$ws.Range("H136").Value = 2692.3684
$ws.Range("I136").Value = 1929.5834
$ws.Range("K136").Value = 5788.7502
$ws.Range("M136").Value = -3238.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 251.16667
$ws.Range("J4").Value = 89
$ws.Range("L4").Value = 89
$ws.Range("N4").Value = -319
$ws.Range("H86").Value = 5281.4
$ws.Range("J86").Value = 5001.75
$ws.Range("L86").Value = 5001.75
$ws.Range("N86").Value = -7247.75
$ws.Range("H89").Value = 5281.4
$ws.Range("J89").Value = 5001.75
$ws.Range("L89").Value = 25008.75
$ws.Range("N89").Value = -36240.75
$ws.Range("H99").Value = 2207.182
$ws.Range("I99").Value = 1880
$ws.Range("K99").Value = 1880
$ws.Range("M99").Value = -382

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3721.2632
$ws.Range("J31").Value = 6137.769
$ws.Range("L31").Value = 6137.769
$ws.Range("N31").Value = -6727.769
$ws.Range("H34").Value = 3721.2632
$ws.Range("J34").Value = 6137.769
$ws.Range("L34").Value = 6137.769
$ws.Range("N34").Value = -6541.769
$ws.Range("H58").Value = 3731.3
$ws.Range("I58").Value = 1052.1666
$ws.Range("K58").Value = 1052.1666
$ws.Range("M58").Value = -849.1666
$ws.Range("H93").Value = 9785.333000000001
$ws.Range("I93").Value = 11342.4
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 11342.4
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = -9470.4
$ws.Range("N93").Value = -5744
$ws.Range("H136").Value = 3731.3
$ws.Range("I136").Value = 1052.1666
$ws.Range("K136").Value = 3156.4998
$ws.Range("M136").Value = -606.4998000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 233.33333
$ws.Range("J33").Value = 100
$ws.Range("L33").Value = 600
$ws.Range("N33").Value = -1166
$ws.Range("H44").Value = 2685.9375
$ws.Range("I44").Value = 259.375
$ws.Range("J44").Value = 5112.5
$ws.Range("K44").Value = 778.125
$ws.Range("L44").Value = 15337.5
$ws.Range("M44").Value = -380.125
$ws.Range("N44").Value = -16133.5
$ws.Range("H46").Value = 1294.625
$ws.Range("J46").Value = 1379.7142
$ws.Range("L46").Value = 4139.142599999999
$ws.Range("N46").Value = -4321.142599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4510.125
$ws.Range("I113").Value = 2328
$ws.Range("J113").Value = 19785
$ws.Range("K113").Value = 2328
$ws.Range("L113").Value = 19785
$ws.Range("M113").Value = -158
$ws.Range("N113").Value = -24125
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3464
$ws.Range("I7").Value = 2883.7144
$ws.Range("K7").Value = 2883.7144
$ws.Range("M7").Value = -2771.7144
$ws.Range("H46").Value = 5065.5
$ws.Range("I46").Value = 3541.5715
$ws.Range("J46").Value = 7199
$ws.Range("K46").Value = 3541.5715
$ws.Range("L46").Value = 7199
$ws.Range("M46").Value = -3353.5715
$ws.Range("N46").Value = -7575
$ws.Range("H68").Value = 5916.5
$ws.Range("I68").Value = 5928.4
$ws.Range("J68").Value = 5896.6665
$ws.Range("K68").Value = 5928.4
$ws.Range("L68").Value = 5896.6665
$ws.Range("M68").Value = -5179.4
$ws.Range("N68").Value = -7394.6665
$ws.Range("H71").Value = 5916.5
$ws.Range("I71").Value = 5928.4
$ws.Range("J71").Value = 5896.6665
$ws.Range("K71").Value = 29642
$ws.Range("L71").Value = 29483.3325
$ws.Range("M71").Value = -25898
$ws.Range("N71").Value = -36971.3325
$ws.Range("H122").Value = 13820.5
$ws.Range("J122").Value = 10839
$ws.Range("L122").Value = 32517
$ws.Range("N122").Value = -37417
$ws.Range("H126").Value = 3464
$ws.Range("I126").Value = 2883.7144
$ws.Range("K126").Value = 8651.143199999999
$ws.Range("M126").Value = -6181.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 1103
$ws.Range("I4").Value = 1750
$ws.Range("J4").Value = 132.5
$ws.Range("K4").Value = 1750
$ws.Range("L4").Value = 132.5
$ws.Range("M4").Value = -1637
$ws.Range("N4").Value = -358.5
$ws.Range("H62").Value = 6177.143
$ws.Range("I62").Value = 5868.8
$ws.Range("K62").Value = 5868.8
$ws.Range("M62").Value = -5244.8
$ws.Range("H65").Value = 6177.143
$ws.Range("I65").Value = 5868.8
$ws.Range("K65").Value = 29344
$ws.Range("M65").Value = -26224
$ws.Range("H126").Value = 2175.5
$ws.Range("I126").Value = 2300.6667
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 6902.000100000001
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -4432.000100000001
$ws.Range("N126").Value = -10340
$ws.Range("H132").Value = 5887.3447
$ws.Range("I132").Value = 5552.778
$ws.Range("K132").Value = 16658.334
$ws.Range("M132").Value = -14128.334
